# "added 4wk low sales check"
# Refreshes the per-week forecast numbers (and the metrics that are
# derived from them) after the forecasting logic learned to treat a
# trailing 4-week run of low sales specially.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Columns on "Forecast Comparison":
#   D = MyForecast, H = Inventory Coverage, I = Stockout Risk,
#   J = Reorder Urgency, L = Seasonality Index
# Row -> [D, H, I, J, L]  ($null = leave untouched)
$rows = @{
    2  = @(9,  13.3,  $null,  $null,    1.05)
    3  = @(9,  11.8,  $null,  $null,    0.96)
    4  = @(10, 10.58, $null,  $null,    0.92)
    5  = @(10, 9.58,  $null,  $null,    0.91)
    6  = @(9,  8.76,  $null,  $null,    1.11)
    7  = @(10, 7.6,   $null,  $null,    0.95)
    8  = @(10, 6.6,   $null,  $null,    0.93)
    9  = @(10, 5.49,  $null,  $null,    1.19)
    10 = @(10, 4.4,   $null,  $null,    0.9399999999999999)
    11 = @(10, 3.4,   $null,  "Normal", 0.96)
    12 = @(10, 2.4,   "Low",  "Normal", 0.84)
    13 = @(10, 1.38,  "Low",  "Normal", 1.14)
    14 = @(10, 0.37,  $null,  $null,    0.88)
    15 = @(10, $null, $null,  $null,    0.99)
    16 = @(10, $null, $null,  $null,    1.04)
    17 = @(10, $null, $null,  $null,    1.13)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]

    $d = $vals[0]
    if ($null -ne $d) { $wsForecast.Cells.Item($r, 4).Value = $d }

    $h = $vals[1]
    if ($null -ne $h) { $wsForecast.Cells.Item($r, 8).Value = $h }

    $i = $vals[2]
    if ($null -ne $i) { $wsForecast.Cells.Item($r, 9).Value = $i }

    $j = $vals[3]
    if ($null -ne $j) { $wsForecast.Cells.Item($r, 10).Value = $j }

    $l = $vals[4]
    if ($null -ne $l) { $wsForecast.Cells.Item($r, 12).Value = $l }
}

# "Summary" sheet metrics recomputed from the refreshed forecast column.
# Column B stores these as text, so force a text format before writing
# the numeric-looking strings (keeps the cell type consistent with the
# other rows on this sheet, e.g. "169 units", "N/A", ...).
$summaryUpdates = @{
    9  = "164"   # Total Forecast (16 Weeks)
    10 = "79"    # Total Forecast (8 Weeks)
    11 = "39"    # Total Forecast (4 Weeks)
    12 = "11"    # Max Forecast
    14 = "9"     # Min Forecast
}

foreach ($r in $summaryUpdates.Keys) {
    $cell = $wsSummary.Cells.Item($r, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$r]
}
